$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - cardholder first name
$ws.Range("C2").Value = "Hartmut"

# Row 3 - card number and surname (keep card number as text, not a number)
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5 - opening balance date
$ws.Range("D5").Value = "KONTOSTAND AM 11.01.2024"

# Row 6 - transaction 1
$ws.Range("B6").Value = "14.01."
$ws.Range("C6").Value = "15.01."
$ws.Range("D6").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 24707432"
$ws.Range("E6").Value = "85,97-"

# Row 7 - transaction 2
$ws.Range("B7").Value = "18.01."
$ws.Range("C7").Value = "19.01."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU KBTCTF"
$ws.Range("E7").Value = "74,37-"

# Row 8 - transaction 3 (amount only changes)
$ws.Range("E8").Value = "48,71-"

# Row 9 - new transaction 4 (previously blank row)
$ws.Range("B9").Value = "25.01."
$ws.Range("C9").Value = "26.01."
$ws.Range("D9").Value = "AMAZON.DE MKTPLC EU XDBTFK"
$ws.Range("E9").HorizontalAlignment = $ws.Range("E8").HorizontalAlignment
$ws.Range("E9").VerticalAlignment = $ws.Range("E8").VerticalAlignment
$ws.Range("E9").WrapText = $ws.Range("E8").WrapText
$ws.Range("E9").Font.Name = $ws.Range("E8").Font.Name
$ws.Range("E9").Font.Size = $ws.Range("E8").Font.Size
$ws.Range("E9").Value = "214,06-"

# Row 12 - closing balance date and amount
$ws.Range("D12").Value = "KONTOSTAND AM 29.01.2024"
$ws.Range("E12").Value = "423,11-"

# Row 13 - next billing date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 08.02.2024"
